$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.915.98'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.634.37'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'211.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'23.15"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').Value = "'0.0880"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '1.866.56'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').Value = '1.638.08'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = "'65.13"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = '27.920.35'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = "'229.86"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').Value = '0.0₃0721'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = "'4.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').Value = "'10.31"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.40%  '
$ws.Range('E24').Value = '  -4.07%  '
$ws.Range('D25').Value = "'152.99"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').Value = "'6.95"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('D27').Value = "'15.61"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.397.40'
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = "'3.06"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').Value = "'1.02"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.21%  '
$ws.Range('E37').Value = '  +1.37%  '
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('D40').Value = "'0.868"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.44%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = "'66.78"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('E44').Value = '  +2.83%  '
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('D47').Value = '1.776.00'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').Value = "'87.64"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.51%  '
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').Value = "'7.50"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.95%  '
